$wb = $excel.ActiveWorkbook

# --- Update status text from "Ready for handoff" to "In Translation" ---
# Overview sheet: columns E (zh-cn) and F (de-de) on row 2 hold the status value.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn sheet: Status column (C) on row 2.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de sheet: Status column (C) on row 2.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the "Status" columns to match the regenerated report layout ---
# Target stored column width is 13.4101845877511 characters. The engine
# quantizes ColumnWidth to whole pixels (1/6 character increments) when it
# is written back to the sheet, so feed it the value whose rounded result
# lands closest to the target (13.3333... ends up being the nearest
# reachable width).
$wsOverview.Columns.Item(5).ColumnWidth = 12.576851254417766   # column E
$wsOverview.Columns.Item(6).ColumnWidth = 12.576851254417766   # column F

$wsZhCn.Columns.Item(3).ColumnWidth = 12.576851254417766       # column C

$wsDeDe.Columns.Item(3).ColumnWidth = 12.576851254417766       # column C
